$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Priming pass: register brand-new shared strings in the desired order ---
# (OffsetF, OffsetA, RD Single, TD Single, then 1Pair-B - matches how the
#  published workbook grew its shared-string table for this commit)
$ws.Range("B7").Value = "OffsetF"
$ws.Range("B8").Value = "OffsetA"
$ws.Range("B9").Value = "RD Single"
$ws.Range("B10").Value = "TD Single"
$ws.Range("L2").Value = "1Pair-B"

# --- Row 1: column index header (0-18), now extends through T1 ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18

# --- Row 2: column labels (HKL row / variable names) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# --- Rows 3-11: one row per method, averaged-intensity ratios ---

# Row 3: Equal Angle
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Equal Angle"
$ws.Range("C3").Value = 1.112773775216138
$ws.Range("D3").Value = 0.8828386167146974
$ws.Range("E3").Value = 0.9821829971181556
$ws.Range("F3").Value = 1.112773775216138
$ws.Range("G3").Value = 0.943393371757925
$ws.Range("H3").Value = 0.970814121037464
$ws.Range("I3").Value = 1.014466858789625
$ws.Range("J3").Value = 0.8828386167146974
$ws.Range("K3").Value = 1.112773775216138
$ws.Range("L3").Value = 0.9821829971181556
$ws.Range("M3").Value = 0.9325108069164265
$ws.Range("N3").Value = 0.9325108069164265
$ws.Range("O3").Value = 0.9361383285302592
$ws.Range("P3").Value = 0.9925984630163306
$ws.Range("Q3").Value = 0.9925984630163304
$ws.Range("R3").Value = 1.022642291066282
$ws.Range("S3").Value = 1.022642291066282
$ws.Range("T3").Value = 0.9844116234390009

# Row 4: CLR
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "CLR"
$ws.Range("C4").Value = 0.9880950006714863
$ws.Range("D4").Value = 0.9971677660083746
$ws.Range("E4").Value = 0.997460550650782
$ws.Range("F4").Value = 0.9880950006714863
$ws.Range("G4").Value = 0.9907443502281067
$ws.Range("H4").Value = 1.00927673239234
$ws.Range("I4").Value = 0.9959477656291794
$ws.Range("J4").Value = 0.9971677660083746
$ws.Range("K4").Value = 0.9880950006714863
$ws.Range("L4").Value = 0.997460550650782
$ws.Range("M4").Value = 0.9973141583295784
$ws.Range("N4").Value = 0.9973141583295784
$ws.Range("O4").Value = 0.9951242222957545
$ws.Range("P4").Value = 0.9942411057768811
$ws.Range("Q4").Value = 0.9942411057768811
$ws.Range("R4").Value = 0.9927045795005324
$ws.Range("S4").Value = 0.9927045795005324
$ws.Range("T4").Value = 0.9964486942633783

# Row 5: BT8Hex
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "BT8Hex"
$ws.Range("C5").Value = 0.9840075537274647
$ws.Range("D5").Value = 0.9942423476527059
$ws.Range("E5").Value = 1.002020827077919
$ws.Range("F5").Value = 0.9840075537274647
$ws.Range("G5").Value = 0.9848882280614535
$ws.Range("H5").Value = 1.026081241052793
$ws.Range("I5").Value = 0.998139362669812
$ws.Range("J5").Value = 0.9942423476527059
$ws.Range("K5").Value = 0.9840075537274647
$ws.Range("L5").Value = 1.002020827077919
$ws.Range("M5").Value = 0.9981315873653125
$ws.Range("N5").Value = 0.9981315873653125
$ws.Range("O5").Value = 0.9937171342640262
$ws.Range("P5").Value = 0.9934235761526966
$ws.Range("Q5").Value = 0.9934235761526966
$ws.Range("R5").Value = 0.9910695705463887
$ws.Range("S5").Value = 0.9910695705463887
$ws.Range("T5").Value = 0.9982299267070246

# Row 6: Spiral
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Spiral"
$ws.Range("C6").Value = 0.9971504836593841
$ws.Range("D6").Value = 0.9928408625209156
$ws.Range("E6").Value = 0.9938793134274173
$ws.Range("F6").Value = 0.9971504836593841
$ws.Range("G6").Value = 0.9926953675046895
$ws.Range("H6").Value = 0.9944493701623758
$ws.Range("I6").Value = 0.9958446873361021
$ws.Range("J6").Value = 0.9928408625209156
$ws.Range("K6").Value = 0.9971504836593841
$ws.Range("L6").Value = 0.9938793134274173
$ws.Range("M6").Value = 0.9933600879741664
$ws.Range("N6").Value = 0.9933600879741664
$ws.Range("O6").Value = 0.9931385144843409
$ws.Range("P6").Value = 0.9946235532025723
$ws.Range("Q6").Value = 0.9946235532025725
$ws.Range("R6").Value = 0.9952552858167754
$ws.Range("S6").Value = 0.9952552858167754
$ws.Range("T6").Value = 0.9944766807684807

# Row 7: OffsetF
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "OffsetF"
$ws.Range("C7").Value = 0.5197175306851636
$ws.Range("D7").Value = 2.567408096928813
$ws.Range("E7").Value = 0.6615123080245794
$ws.Range("F7").Value = 0.5197175306851636
$ws.Range("G7").Value = 1.791940907618374
$ws.Range("H7").Value = 0.2152412730765806
$ws.Range("I7").Value = 0.6150759560616786
$ws.Range("J7").Value = 2.567408096928813
$ws.Range("K7").Value = 0.5197175306851636
$ws.Range("L7").Value = 0.6615123080245794
$ws.Range("M7").Value = 1.614460202476696
$ws.Range("N7").Value = 1.614460202476696
$ws.Range("O7").Value = 1.673620437523922
$ws.Range("P7").Value = 1.249545978546185
$ws.Range("Q7").Value = 1.249545978546185
$ws.Range("R7").Value = 1.06708886658093
$ws.Range("S7").Value = 1.06708886658093
$ws.Range("T7").Value = 1.061816012065865

# Row 8: OffsetA
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "OffsetA"
$ws.Range("C8").Value = 0.9217430670707433
$ws.Range("D8").Value = 1.282746925950709
$ws.Range("E8").Value = 0.9058280755210715
$ws.Range("F8").Value = 0.9217430670707433
$ws.Range("G8").Value = 1.217755207895157
$ws.Range("H8").Value = 0.637312924631213
$ws.Range("I8").Value = 0.9089803550235943
$ws.Range("J8").Value = 1.282746925950709
$ws.Range("K8").Value = 0.9217430670707433
$ws.Range("L8").Value = 0.9058280755210715
$ws.Range("M8").Value = 1.09428750073589
$ws.Range("N8").Value = 1.09428750073589
$ws.Range("O8").Value = 1.135443403122313
$ws.Range("P8").Value = 1.036772689514175
$ws.Range("Q8").Value = 1.036772689514175
$ws.Range("R8").Value = 1.008015283903317
$ws.Range("S8").Value = 1.008015283903317
$ws.Range("T8").Value = 0.9790610926820814

# Row 9: RD Single
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "RD Single"
$ws.Range("C9").Value = 1.97
$ws.Range("D9").Value = 0.21
$ws.Range("E9").Value = 0.86
$ws.Range("F9").Value = 1.97
$ws.Range("G9").Value = 0.64
$ws.Range("H9").Value = 0.67
$ws.Range("I9").Value = 1.11
$ws.Range("J9").Value = 0.21
$ws.Range("K9").Value = 1.97
$ws.Range("L9").Value = 0.86
$ws.Range("M9").Value = 0.535
$ws.Range("N9").Value = 0.535
$ws.Range("O9").Value = 0.57
$ws.Range("P9").Value = 1.013333333333333
$ws.Range("Q9").Value = 1.013333333333333
$ws.Range("R9").Value = 1.2525
$ws.Range("S9").Value = 1.2525
$ws.Range("T9").Value = 0.9100000000000001

# Row 10: TD Single
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "TD Single"
$ws.Range("C10").Value = 0.21
$ws.Range("D10").Value = 3.94
$ws.Range("E10").Value = 0.46
$ws.Range("F10").Value = 0.21
$ws.Range("G10").Value = 2.23
$ws.Range("H10").Value = 0.07000000000000001
$ws.Range("I10").Value = 0.4
$ws.Range("J10").Value = 3.94
$ws.Range("K10").Value = 0.21
$ws.Range("L10").Value = 0.46
$ws.Range("M10").Value = 2.2
$ws.Range("N10").Value = 2.2
$ws.Range("O10").Value = 2.21
$ws.Range("P10").Value = 1.536666666666667
$ws.Range("Q10").Value = 1.536666666666667
$ws.Range("R10").Value = 1.205
$ws.Range("S10").Value = 1.205
$ws.Range("T10").Value = 1.218333333333333

# Row 11: HexGrid-90degTilt5degRes
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C11").Value = 0.9954624986232038
$ws.Range("D11").Value = 0.993213724390445
$ws.Range("E11").Value = 0.9942589894037912
$ws.Range("F11").Value = 0.9954624986232038
$ws.Range("G11").Value = 0.9934016870969857
$ws.Range("H11").Value = 0.993686845332701
$ws.Range("I11").Value = 0.995686040388355
$ws.Range("J11").Value = 0.993213724390445
$ws.Range("K11").Value = 0.9954624986232038
$ws.Range("L11").Value = 0.9942589894037912
$ws.Range("M11").Value = 0.993736356897118
$ws.Range("N11").Value = 0.993736356897118
$ws.Range("O11").Value = 0.9936248002970739
$ws.Range("P11").Value = 0.99431173747248
$ws.Range("Q11").Value = 0.99431173747248
$ws.Range("R11").Value = 0.994599427760161
$ws.Range("S11").Value = 0.994599427760161
$ws.Range("T11").Value = 0.9942849642059136

# --- Formatting: extend the bold/centered/bordered style (style index 1)
#     from the existing header cells onto the newly added ones ---
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
